# Apply cryptos-list update: refresh Price (D) and Volume(1h) (E) columns.
# Values are written as text (matching the source data's inline-string cells);
# Excel would otherwise auto-coerce numeric-looking strings (e.g. "172.90", "1.00")
# into numbers and silently drop the trailing zero, so we force text via a
# temporary "@" (Text) number format and then clear the format again so the
# cell style is left exactly as it was (no residual style index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '70.814.88'
Set-TextValue $ws.Range("E2") '  +2.78%  '

Set-TextValue $ws.Range("D3") '3.795.66'
Set-TextValue $ws.Range("E3") '  +0.87%  '

Set-TextValue $ws.Range("E4") '  +0.04%  '

Set-TextValue $ws.Range("D5") '702.14'
Set-TextValue $ws.Range("E5") '  +10.18%  '

Set-TextValue $ws.Range("D6") '172.90'
Set-TextValue $ws.Range("E6") '  +4.47%  '

Set-TextValue $ws.Range("D7") '3.796.27'
Set-TextValue $ws.Range("E7") '  +0.90%  '

Set-TextValue $ws.Range("E8") '  +0.03%  '

Set-TextValue $ws.Range("E9") '  +1.01%  '

Set-TextValue $ws.Range("E10") '  +2.61%  '

Set-TextValue $ws.Range("D11") '7.39'
Set-TextValue $ws.Range("E11") '  +6.79%  '

Set-TextValue $ws.Range("D12") '0.460'
Set-TextValue $ws.Range("E12") '  +0.89%  '

Set-TextValue $ws.Range("D13") '0.0000258'
Set-TextValue $ws.Range("E13") '  +8.15%  '

Set-TextValue $ws.Range("D14") '36.37'
Set-TextValue $ws.Range("E14") '  +4.32%  '

Set-TextValue $ws.Range("D15") '4.435.34'
Set-TextValue $ws.Range("E15") '  +0.92%  '

Set-TextValue $ws.Range("D16") '3.797.69'
Set-TextValue $ws.Range("E16") '  +1.06%  '

Set-TextValue $ws.Range("D17") '70.805.91'
Set-TextValue $ws.Range("E17") '  +2.87%  '

Set-TextValue $ws.Range("D18") '17.86'
Set-TextValue $ws.Range("E18") '  +1.09%  '

Set-TextValue $ws.Range("D19") '7.19'
Set-TextValue $ws.Range("E19") '  +2.91%  '

Set-TextValue $ws.Range("E20") '  +0.47%  '

Set-TextValue $ws.Range("E21") '  +16.81%  '

Set-TextValue $ws.Range("D22") '481.81'
Set-TextValue $ws.Range("E22") '  +2.43%  '

Set-TextValue $ws.Range("E23") '  +1.68%  '

Set-TextValue $ws.Range("D24") '84.61'
Set-TextValue $ws.Range("E24") '  +3.68%  '

Set-TextValue $ws.Range("D25") '0.0000144'
Set-TextValue $ws.Range("E25") '  +0.71%  '

Set-TextValue $ws.Range("D26") '12.40'
Set-TextValue $ws.Range("E26") '  +2.05%  '

Set-TextValue $ws.Range("E27") '  +3.49%  '

Set-TextValue $ws.Range("D28") '10.44'
Set-TextValue $ws.Range("E28") '  +4.19%  '

Set-TextValue $ws.Range("D29") '3.945.32'
Set-TextValue $ws.Range("E29") '  +0.89%  '

Set-TextValue $ws.Range("D30") '1.00'
Set-TextValue $ws.Range("E30") '  -0.08%  '

Set-TextValue $ws.Range("D31") '3.15'
Set-TextValue $ws.Range("E31") '  +17.35%  '

Set-TextValue $ws.Range("D32") '7.53'
Set-TextValue $ws.Range("E32") '  +6.12%  '

Set-TextValue $ws.Range("D33") '2.28'
Set-TextValue $ws.Range("E33") '  +0.87%  '

Set-TextValue $ws.Range("E34") '  +3.90%  '

Set-TextValue $ws.Range("E35") '  +4.67%  '

Set-TextValue $ws.Range("D36") '9.24'
Set-TextValue $ws.Range("E36") '  +4.15%  '

Set-TextValue $ws.Range("E37") '  -0.09%  '

Set-TextValue $ws.Range("E38") '  +2.20%  '

Set-TextValue $ws.Range("D39") '3.44'
Set-TextValue $ws.Range("E39") '  +6.29%  '

Set-TextValue $ws.Range("E40") '  +4.68%  '

Set-TextValue $ws.Range("E41") '  +12.53%  '

Set-TextValue $ws.Range("D42") '0.972'
Set-TextValue $ws.Range("E42") '  +1.87%  '

Set-TextValue $ws.Range("D43") '0.000326'
Set-TextValue $ws.Range("E43") '  +22.32%  '

Set-TextValue $ws.Range("D44") '0.999'
Set-TextValue $ws.Range("E44") '  +0.03%  '

Set-TextValue $ws.Range("E45") '  +0.02%  '

Set-TextValue $ws.Range("D46") '162.63'
Set-TextValue $ws.Range("E46") '  +4.64%  '

Set-TextValue $ws.Range("D47") '49.05'
Set-TextValue $ws.Range("E47") '  +3.58%  '

Set-TextValue $ws.Range("D48") '44.79'
Set-TextValue $ws.Range("E48") '  -0.05%  '

Set-TextValue $ws.Range("E49") '  +3.25%  '

Set-TextValue $ws.Range("E50") '  -1.13%  '

Set-TextValue $ws.Range("D51") '8.56'
Set-TextValue $ws.Range("E51") '  +2.56%  '
